$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 964.9583
$ws.Range("I33").Value = 279.95456
$ws.Range("J33").Value = 8500
$ws.Range("K33").Value = 279.95456
$ws.Range("L33").Value = 8500
$ws.Range("M33").Value = -50.95456000000001
$ws.Range("N33").Value = -8958

$ws.Range("H69").Value = 5961.4287
$ws.Range("J69").Value = 6314.5454
$ws.Range("L69").Value = 18943.6362
$ws.Range("N69").Value = -20691.6362

$ws.Range("H72").Value = 5961.4287
$ws.Range("J72").Value = 6314.5454
$ws.Range("L72").Value = 56830.9086
$ws.Range("N72").Value = -65566.9086

$ws.Range("H116").Value = 2215.5
$ws.Range("I116").Value = 1763.5
$ws.Range("J116").Value = 2328.5
$ws.Range("K116").Value = 1763.5
$ws.Range("L116").Value = 2328.5
$ws.Range("M116").Value = 1678.5
$ws.Range("N116").Value = -9212.5

$ws.Range("H137").Value = 1538.0714
$ws.Range("I137").Value = 1224.6666
$ws.Range("J137").Value = 10000
$ws.Range("K137").Value = 3673.9998
$ws.Range("L137").Value = 30000
$ws.Range("M137").Value = -1123.9998
$ws.Range("N137").Value = -35100

$ws.Range("H138").Value = 6508.137
$ws.Range("I138").Value = 1255.7142
$ws.Range("J138").Value = 17997.812
$ws.Range("K138").Value = 3767.1426
$ws.Range("L138").Value = 53993.436
$ws.Range("M138").Value = 1372.8574
$ws.Range("N138").Value = -64273.436

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H61").Value = 1506.4736
$ws.Range("I61").Value = 1121.5333
$ws.Range("J61").Value = 2950
$ws.Range("K61").Value = 1121.5333
$ws.Range("L61").Value = 2950
$ws.Range("M61").Value = -909.5333000000001
$ws.Range("N61").Value = -3374

$ws.Range("H74").Value = 3687.4285
$ws.Range("I74").Value = 10012
$ws.Range("K74").Value = 10012
$ws.Range("M74").Value = -9138

$ws.Range("H77").Value = 3687.4285
$ws.Range("I77").Value = 10012
$ws.Range("K77").Value = 50060
$ws.Range("M77").Value = -45692

$ws.Range("H132").Value = 2079.5334
$ws.Range("I132").Value = 2044.3137
$ws.Range("J132").Value = 2279.111
$ws.Range("K132").Value = 6132.9411
$ws.Range("L132").Value = 6837.333
$ws.Range("M132").Value = -3602.9411
$ws.Range("N132").Value = -11897.333

$ws.Range("H136").Value = 1506.4736
$ws.Range("I136").Value = 1121.5333
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 3364.5999
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -814.5999000000002
$ws.Range("N136").Value = -13950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 19900
$ws.Range("J19").Value = 19900
$ws.Range("L19").Value = 19900
$ws.Range("N19").Value = -20246

$ws.Range("H35").Value = 19163.5
$ws.Range("J35").Value = 19163.5
$ws.Range("L35").Value = 19163.5
$ws.Range("N35").Value = -19783.5

$ws.Range("H82").Value = 14173.728
$ws.Range("I82").Value = 2477
$ws.Range("J82").Value = 23921
$ws.Range("K82").Value = 2477
$ws.Range("L82").Value = 23921
$ws.Range("M82").Value = -2094
$ws.Range("N82").Value = -24687

$ws.Range("H85").Value = 14173.728
$ws.Range("I85").Value = 2477
$ws.Range("J85").Value = 23921
$ws.Range("K85").Value = 2477
$ws.Range("L85").Value = 23921
$ws.Range("M85").Value = -1151
$ws.Range("N85").Value = -26573

$ws.Range("H134").Value = 2337.6072
$ws.Range("I134").Value = 2149.875
$ws.Range("J134").Value = 3464
$ws.Range("K134").Value = 6449.625
$ws.Range("L134").Value = 10392
$ws.Range("M134").Value = -3914.625
$ws.Range("N134").Value = -15462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3597.625
$ws.Range("I58").Value = 868.8
$ws.Range("K58").Value = 868.8
$ws.Range("M58").Value = -665.8

$ws.Range("H68").Value = 18308.523
$ws.Range("J68").Value = 18308.523
$ws.Range("L68").Value = 18308.523
$ws.Range("N68").Value = -19806.523

$ws.Range("H71").Value = 18308.523
$ws.Range("J71").Value = 18308.523
$ws.Range("L71").Value = 54925.569
$ws.Range("N71").Value = -62413.569

$ws.Range("H74").Value = 26179.143
$ws.Range("J74").Value = 26179.143
$ws.Range("L74").Value = 26179.143
$ws.Range("N74").Value = -27927.143

$ws.Range("H77").Value = 26179.143
$ws.Range("J77").Value = 26179.143
$ws.Range("L77").Value = 78537.429
$ws.Range("N77").Value = -87273.429

$ws.Range("H134").Value = 803.24445
$ws.Range("I134").Value = 715.14703
$ws.Range("K134").Value = 2145.44109
$ws.Range("M134").Value = 389.5589100000002

$ws.Range("H136").Value = 3597.625
$ws.Range("I136").Value = 868.8
$ws.Range("K136").Value = 2606.4
$ws.Range("M136").Value = -56.39999999999964

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2442.347
$ws.Range("I132").Value = 1505.4062
$ws.Range("J132").Value = 4206
$ws.Range("K132").Value = 4516.2186
$ws.Range("L132").Value = 12618
$ws.Range("M132").Value = -1986.2186
$ws.Range("N132").Value = -17678

$ws.Range("H135").Value = 27921.75
$ws.Range("J135").Value = 27921.75
$ws.Range("L135").Value = 27921.75
$ws.Range("N135").Value = -38061.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3436.6155
$ws.Range("I132").Value = 3165.516
$ws.Range("J132").Value = 4487.125
$ws.Range("K132").Value = 9496.548000000001
$ws.Range("L132").Value = 13461.375
$ws.Range("M132").Value = -6966.548000000001
$ws.Range("N132").Value = -18521.375

$ws.Range("H136").Value = 2084.9443
$ws.Range("I136").Value = 1837.6154
$ws.Range("J136").Value = 2728
$ws.Range("K136").Value = 5512.8462
$ws.Range("L136").Value = 8184
$ws.Range("M136").Value = -2962.8462
$ws.Range("N136").Value = -13284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H126").Value = 2173.8572
$ws.Range("I126").Value = 2076.8
$ws.Range("K126").Value = 6230.400000000001
$ws.Range("M126").Value = -3760.400000000001

$ws.Range("H132").Value = 2054.6724
$ws.Range("I132").Value = 1808.3269
$ws.Range("J132").Value = 4189.6665
$ws.Range("K132").Value = 5424.9807
$ws.Range("L132").Value = 12568.9995
$ws.Range("M132").Value = -2894.9807

Write-Host "Applied all market price updates"
